$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated simulation results (more games simulated) to the transition matrix
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5333333333333333
$ws.Range("S2").Value = 0.2666666666666667
$ws.Range("P4").Value = 1
$ws.Range("F6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.1818181818181818
$ws.Range("Q6").Value = 0.09090909090909091
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.5454545454545454
$ws.Range("Q7").Value = 0.125
$ws.Range("R7").Value = 0.125
$ws.Range("S7").Value = 0.75
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.02777777777777778
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.08333333333333333
$ws.Range("O8").Value = 0.05555555555555555
$ws.Range("Q8").Value = 0.02777777777777778
$ws.Range("R8").Value = 0.08333333333333333
$ws.Range("S8").Value = 0.5555555555555556
$ws.Range("F9").Value = 0.03846153846153846
$ws.Range("Q9").Value = 0.1538461538461539
$ws.Range("R9").Value = 0.1153846153846154
$ws.Range("B10").Value = 0.07792207792207792
$ws.Range("D10").Value = 0.01298701298701299
$ws.Range("F10").Value = 0.03896103896103896
$ws.Range("J10").Value = 0.1038961038961039
$ws.Range("Q10").Value = 0.2207792207792208
$ws.Range("R10").Value = 0.02597402597402598
$ws.Range("S10").Value = 0.5194805194805194
$ws.Range("G11").Value = 0.1538461538461539
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.1923076923076923
$ws.Range("L11").Value = 0.5384615384615384
$ws.Range("S11").Value = 0.03846153846153846
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("S12").Value = 0.1428571428571428
$ws.Range("G13").Value = 1
$ws.Range("F15").Value = 0.06666666666666667
$ws.Range("H15").Value = 0.3333333333333333
$ws.Range("I15").Value = 0.2
$ws.Range("J15").Value = 0.2
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.1333333333333333
$ws.Range("F16").Value = 0.1111111111111111
$ws.Range("H16").Value = 0.1111111111111111
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.4444444444444444
$ws.Range("M16").Value = 0.1111111111111111
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.04
$ws.Range("H17").Value = 0.16
$ws.Range("I17").Value = 0.08
$ws.Range("J17").Value = 0.32
$ws.Range("K17").Value = 0.16
$ws.Range("O17").Value = 0.12
$ws.Range("S17").Value = 0.12
$ws.Range("H18").Value = 0.2727272727272727
$ws.Range("I18").Value = 0.2727272727272727
$ws.Range("J18").Value = 0.1818181818181818
$ws.Range("O18").Value = 0.09090909090909091
$ws.Range("S18").Value = 0.1818181818181818
$ws.Range("F19").Value = 0.00819672131147541
$ws.Range("H19").Value = 0.180327868852459
$ws.Range("I19").Value = 0.139344262295082
$ws.Range("J19").Value = 0.3360655737704918
$ws.Range("K19").Value = 0.1311475409836066
$ws.Range("M19").Value = 0.00819672131147541
